# The "Recorded By" column (G) lists the sources that recorded each
# session's attendance. Previously entries combining "System" and the
# coordinator's email were ordered as "System, <email>". Flip the order
# so the email comes first: "<email>, System".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns.Item(7).Replace("System, dnasr281@gmail.com", "dnasr281@gmail.com, System", 1)
